# Append the latest "Market Open" price-tick rows scraped at 2025-05-12
# 19:xx UTC to each ticker sheet.

$wb = $excel.ActiveWorkbook

function Set-TextCell($sheet, $row, $col, $val) {
    # Assigning a date-/currency-looking string straight to .Value lets
    # Excel "smart" parse it into a real date/number (and stamp a style on
    # the cell). Forcing text format first keeps it a literal string; then
    # ClearFormats drops the now-unneeded number-format override so the
    # cell matches its neighbours (no explicit style).
    $cell = $sheet.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

function Add-TickRow($sheet, $row, $date, $time, $priceText, $priceFloat, $note) {
    Set-TextCell $sheet $row 1 $date
    Set-TextCell $sheet $row 2 $time
    Set-TextCell $sheet $row 3 $priceText
    $sheet.Cells.Item($row, 4).Value = $priceFloat
    Set-TextCell $sheet $row 5 $note
}

$tesla = $wb.Worksheets.Item("Tesla")
Add-TickRow $tesla 20 "2025-05-12" "19:08:39" '$319.01' 319.01 "Market Open"
Add-TickRow $tesla 21 "2025-05-12" "19:20:02" '$318.16' 318.16 "Market Open"

$apple = $wb.Worksheets.Item("Apple")
Add-TickRow $apple 17 "2025-05-12" "19:10:04" '$210.86' 210.86 "Market Open"
Add-TickRow $apple 18 "2025-05-12" "19:20:09" '$210.93' 210.93 "Market Open"

$nvidia = $wb.Worksheets.Item("Nvidia")
Add-TickRow $nvidia 17 "2025-05-12" "19:07:58" '$122.79' 122.79 "Market Open"
Add-TickRow $nvidia 18 "2025-05-12" "19:20:16" '$122.80' 122.8 "Market Open"

$manchester = $wb.Worksheets.Item("Manchester")
Add-TickRow $manchester 18 "2025-05-12" "19:19:36" '$14.57' 14.57 "Market Open"
Add-TickRow $manchester 19 "2025-05-12" "19:20:38" '$14.57' 14.57 "Market Open"
Add-TickRow $manchester 20 "2025-05-12" "19:20:53" '$14.57' 14.57 "Market Open"

$google = $wb.Worksheets.Item("Google")
Add-TickRow $google 18 "2025-05-12" "18:52:45" '$158.36' 158.36 "Market Open"
Add-TickRow $google 19 "2025-05-12" "19:20:23" '$158.38' 158.38 "Market Open"

$nike = $wb.Worksheets.Item("Nike")
Add-TickRow $nike 15 "2025-05-12" "19:20:30" '$62.44' 62.44 "Market Open"
Add-TickRow $nike 16 "2025-05-12" "19:21:05" '$62.45' 62.45 "Market Open"
